# Mẫu 26 – remove the leftover "vnpt.SiteAddress" merge-field placeholder
# text that follows "Địa chỉ: " in the bullet list near the top of the
# document. The whole run carrying that text is deleted (not just its
# text), which Find/Replace achieves by replacing the found text with an
# empty string and Word collapsing the now-empty run.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$found = $find.Execute(
    "vnpt.SiteAddress",  # FindText
    $false,              # MatchCase
    $false,              # MatchWholeWord
    $false,              # MatchWildcards
    $false,              # MatchSoundsLike
    $false,              # MatchAllWordForms
    $true,               # Forward
    1,                   # Wrap (wdFindContinue)
    $false,              # Format
    "",                  # ReplaceWith
    2                    # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Could not find 'vnpt.SiteAddress' text to remove"
}
